$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "q1" in C1, centered
$ws.Range("C1").Value = "q1"
$ws.Range("C1").HorizontalAlignment = -4108  # xlCenter

# Add value 4 in C2 (AGUDELO MORENO RENATO DIRNEY row)
$ws.Range("C2").Value = 4

# Update selection to C1 like in the diff
$ws.Range("C1").Select()
